# Update the division-problem worksheet numbers to match the new output.
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "803÷9=" "326÷4="
Replace-Text "873÷5=" "680÷4="
Replace-Text "672÷3=" "189÷7="
Replace-Text "968÷2=" "490÷8="
Replace-Text "210÷3=" "872÷2="

Replace-Text "620÷7=" "690÷8="
Replace-Text "366÷8=" "759÷4="
Replace-Text "601÷5=" "159÷2="
Replace-Text "598÷3=" "678÷6="
Replace-Text "693÷8=" "118÷2="

Replace-Text "332÷5=" "864÷5="
Replace-Text "149÷2=" "819÷9="
Replace-Text "176÷8=" "667÷7="
Replace-Text "523÷6=" "498÷6="
Replace-Text "235÷2=" "816÷4="

Replace-Text "588÷3=" "302÷5="
Replace-Text "889÷7=" "185÷7="
Replace-Text "101÷6=" "891÷4="
Replace-Text "644÷2=" "974÷8="
Replace-Text "898÷9=" "493÷6="

# The fifth row has "479÷5=" twice (cols 1 and 2), which diverge into two
# different values, so address them positionally via the table rather than
# a text-wide Find/Replace.
$t = $d.Tables(1)
$t.Rows(17).Cells(1).Range.Text = "678÷6="
$t.Rows(17).Cells(2).Range.Text = "485÷3="

Replace-Text "832÷9=" "610÷2="
Replace-Text "257÷9=" "823÷3="
Replace-Text "180÷9=" "115÷5="
